# Insert a new data row before the current row 200 (pushing all rows
# 200-305 down to 201-306, so the final sheet spans A1:R306), then
# populate the newly inserted row with the new "Jengibre" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 200 (and everything below it) down by one row.
$ws.Rows.Item(200).Insert()

# Fill in the new row 200 with the new weekly record.
$ws.Range("A200").Value = 10
$ws.Range("B200").Value = "Vega Modelo de Temuco"
$ws.Range("C200").Value = "La Araucanía"
$ws.Range("D200").Value = 45089
$ws.Range("E200").Value = 9
$ws.Range("F200").Value = 100114007
$ws.Range("G200").Value = "Jengibre"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 55
$ws.Range("K200").Value = 24000
$ws.Range("L200").Value = 24000
$ws.Range("M200").Value = 24000
$ws.Range("N200").Value = "$/caja 13 kilos"
$ws.Range("O200").Value = "Perú"
$ws.Range("P200").Value = 1846
$ws.Range("Q200").Value = 13
$ws.Range("R200").Value = "Hortaliza"
